$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (prices & 1h volume change)
# Columns B (Coin) and C (Link) stay text naturally; D (Price) and E (Volume)
# must stay text too, so force text format before assigning to avoid Excel
# auto-converting numeric-looking strings like "1.00" into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.331.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.642.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.82"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.84%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.11%  "

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.172"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +16.14%  "

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.638.64"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.165"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.13%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.92%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +10.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.084.20"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "72.141.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.33%  "

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.71"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.632.98"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.11%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +7.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "384.20"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.40%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.42%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +21.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.17"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.46"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.59%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +11.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.775.68"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.19%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0973"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +11.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "550.38"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.10"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.67%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.94%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.62"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.31"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.39%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.15"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.76%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +10.87%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.27%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +14.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.15%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.45"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.13%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.540"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.50%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +9.48%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.76%  "
